# Added data for Pb ads experiments
# Populates column F (pH) for rows 2-55 on Sheet1 with the newly measured
# values. Two of the readings (rows 7 and 54) are missing/invalid and are
# recorded as the literal text "nan" (same convention already used
# elsewhere in the sheet, e.g. column D), everything else is numeric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fValues = [ordered]@{
    2  = 5.15
    3  = 5.306
    4  = 5.3
    5  = 5.782
    6  = 5.904
    7  = "nan"
    8  = 5.15
    9  = 6.029
    10 = 5.756
    11 = 6.08
    12 = 6.084
    13 = 5.914
    14 = 5.15
    15 = 5.757
    16 = 5.71
    17 = 5.98
    18 = 5.829
    19 = 5.283
    20 = 5.01
    21 = 5.06
    22 = 5.101
    23 = 4.594
    24 = 4.932
    25 = 5.291
    26 = 5.01
    27 = 5.229
    28 = 4.834
    29 = 5.279
    30 = 5.378
    31 = 5
    32 = 5.01
    33 = 5.321
    34 = 5.109
    35 = 5.188
    36 = 5.136
    37 = 4.942
    38 = 5.05
    39 = 4.782
    40 = 4.946
    41 = 5.253
    42 = 5.199
    43 = 4.954
    44 = 5.05
    45 = 5.068
    46 = 5.132
    47 = 5.146
    48 = 5.27
    49 = 4.873
    50 = 5.05
    51 = 4.84
    52 = 5.315
    53 = 5.109
    54 = "nan"
    55 = 4.972
}

foreach ($row in $fValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $fValues[$row]
}

# Leave the selection where the author finished entering the new data
# (next empty row below the new F column values).
$ws.Range("F56").Select() | Out-Null

